# Se corrige el error de la automatización al seleccionar el número de habitaciones
# The "Tarifa" sheet stores the number of days (NUMERODEDIAS) used by the
# data-driven automation test in cell B2. It was incorrectly set to 14 and
# must be corrected to 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tarifa")

$ws.Range("B2").Value = 13
